# Update TPM-derived values in the LR-pairs sheet (Fgf8-Fgfr1) with newly
# computed results from the updated scripts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 2.294987
$ws.Range("N2").Value = 6.884961000000001
$ws.Range("O2").Value = 0.0158275801650097
$ws.Range("P2").Value = 0.0158275801650097
$ws.Range("Q2").Value = 0.2955468958686667
$ws.Range("R2").Value = 2.659922062818
$ws.Range("S2").Value = 0.0158275801650097
$ws.Range("T2").Value = 0.0158275801650097

# Row 3
$ws.Range("O3").Value = 0.769602070219672
$ws.Range("P3").Value = 0.7696020702196722
$ws.Range("S3").Value = 0.769602070219672
$ws.Range("T3").Value = 0.7696020702196722

# Row 4
$ws.Range("M4").Value = 31.11253633333333
$ws.Range("N4").Value = 93.337609
$ws.Range("O4").Value = 0.2145703496153182
$ws.Range("P4").Value = 0.2145703496153182
$ws.Range("R4").Value = 36.059865185842
$ws.Range("S4").Value = 0.2145703496153182
$ws.Range("T4").Value = 0.2145703496153182
